$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B column values for rows 5 through 12
$ws.Range("B5").Value = 403940
$ws.Range("B6").Value = 424940
$ws.Range("B7").Value = 434940
$ws.Range("B8").Value = 444940
$ws.Range("B9").Value = 444940
$ws.Range("B10").Value = 426435
$ws.Range("B11").Value = 408467
$ws.Range("B12").Value = 400267

# Update the selected range shown when the workbook is reopened
$ws.Range("A1:A20").Select()
